$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Teacher ("guru") rows: split the old "guru1/guru2/guru3" single name
# into a shared first-name "guru" (col A) + distinct last-name suffix (col B),
# and rebuild the e-mail accordingly (col C). Hyperlinks / styles on C2:C4
# are left untouched (only their displayed text changes).
$ws.Range("A2").Value = "guru"
$ws.Range("B2").Value = "bin"
$ws.Range("C2").Value = "gurubin@tes.com"

$ws.Range("A3").Value = "guru"
$ws.Range("B3").Value = "mtk"
$ws.Range("C3").Value = "gurumtk@tes.com"

$ws.Range("A4").Value = "guru"
$ws.Range("B4").Value = "big"
$ws.Range("C4").Value = "gurubig@tes.com"

# --- Student ("siswa") rows 5-13: duplicate first name into new column B
# (keeps the existing formula-driven e-mail in column C as-is).
For ($r = 5; $r -le 13; $r++) {
    $ws.Range("B$r").Value = $ws.Range("A$r").Value2
}

# --- Drop students siswa10..siswa20 (previously rows 14-24) entirely.
# Use ClearContents (not a row delete) so the full-column dataValidation
# ranges below (D1:D1048576 etc.) are not shifted/renumbered.
$ws.Range("A14:F24").ClearContents()

# --- Restore the selection to match the edited range.
[void]$ws.Range("F6:F13").Select()
